$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.666448
$ws.Range("H2").Value = 3.332896
$ws.Range("I2").Value = 0.08698774157534103
$ws.Range("J2").Value = 0.0653244391585683
$ws.Range("M2").Value = 35.18694
$ws.Range("N2").Value = 70.37388
$ws.Range("O2").Value = 0.1785572969025014
$ws.Range("P2").Value = 0.1282643028201015
$ws.Range("Q2").Value = 58.63720578912
$ws.Range("R2").Value = 234.54882315648
$ws.Range("S2").Value = 0.01553229599934623
$ws.Range("T2").Value = 0.0083787936457879

$ws.Range("G3").Value = 1.666448
$ws.Range("H3").Value = 3.332896
$ws.Range("I3").Value = 0.08698774157534103
$ws.Range("J3").Value = 0.0653244391585683
$ws.Range("O3").Value = 0.02719713854783473
$ws.Range("P3").Value = 0.02930506404712654
$ws.Range("Q3").Value = 8.931386381679999
$ws.Range("R3").Value = 53.58831829008
$ws.Range("S3").Value = 0.002365817659587793
$ws.Range("T3").Value = 0.001914336873384465

$ws.Range("G4").Value = 1.666448
$ws.Range("H4").Value = 3.332896
$ws.Range("I4").Value = 0.08698774157534103
$ws.Range("J4").Value = 0.0653244391585683
$ws.Range("M4").Value = 70.59161999999999
$ws.Range("N4").Value = 211.77486
$ws.Range("O4").Value = 0.3582195226742806
$ws.Range("P4").Value = 0.3859834752997077
$ws.Range("Q4").Value = 117.63726396576
$ws.Range("R4").Value = 705.8235837945599
$ws.Range("S4").Value = 0.03116070726563234
$ws.Range("T4").Value = 0.02521415404842851

$ws.Range("G5").Value = 1.666448
$ws.Range("H5").Value = 3.332896
$ws.Range("I5").Value = 0.08698774157534103
$ws.Range("J5").Value = 0.0653244391585683
$ws.Range("M5").Value = 7.337415
$ws.Range("N5").Value = 14.67483
$ws.Range("O5").Value = 0.03723395636710288
$ws.Range("P5").Value = 0.0267465263952124
$ws.Range("Q5").Value = 12.22742055192
$ws.Range("R5").Value = 48.90968220768
$ws.Range("S5").Value = 0.003238897774289069
$ws.Range("T5").Value = 0.001747201836207093

$ws.Range("G6").Value = 1.666448
$ws.Range("H6").Value = 3.332896
$ws.Range("I6").Value = 0.08698774157534103
$ws.Range("J6").Value = 0.0653244391585683
$ws.Range("M6").Value = 15.36873766666667
$ws.Range("N6").Value = 46.106213
$ws.Range("O6").Value = 0.07798917025929666
$ws.Range("P6").Value = 0.0840337532351506
$ws.Range("Q6").Value = 25.61120214714133
$ws.Range("R6").Value = 153.667212882848
$ws.Range("S6").Value = 0.006784101788190971
$ws.Range("T6").Value = 0.005489457800475737

$ws.Range("G7").Value = 1.666448
$ws.Range("H7").Value = 3.332896
$ws.Range("I7").Value = 0.08698774157534103
$ws.Range("J7").Value = 0.0653244391585683
$ws.Range("M7").Value = 63.21821133333334
$ws.Range("N7").Value = 189.654634
$ws.Range("O7").Value = 0.3208029152489838
$ws.Range("P7").Value = 0.3456668782027013
$ws.Range("Q7").Value = 105.3498618400107
$ws.Range("R7").Value = 632.099171040064
$ws.Range("S7").Value = 0.02790592108829463
$ws.Range("T7").Value = 0.0225804949542846

$ws.Range("I8").Value = 0.1669502665149541
$ws.Range("J8").Value = 0.1880598173367416
$ws.Range("M8").Value = 35.18694
$ws.Range("N8").Value = 70.37388
$ws.Range("O8").Value = 0.1785572969025014
$ws.Range("P8").Value = 0.1282643028201015
$ws.Range("Q8").Value = 112.53881244528
$ws.Range("R8").Value = 675.2328746716801
$ws.Range("S8").Value = 0.02981018830606238
$ws.Range("T8").Value = 0.0241213613591728

$ws.Range("I9").Value = 0.1669502665149541
$ws.Range("J9").Value = 0.1880598173367416
$ws.Range("O9").Value = 0.02719713854783473
$ws.Range("P9").Value = 0.02930506404712654
$ws.Range("S9").Value = 0.004540569529005139
$ws.Range("T9").Value = 0.005511104991744132

$ws.Range("I10").Value = 0.1669502665149541
$ws.Range("J10").Value = 0.1880598173367416
$ws.Range("M10").Value = 70.59161999999999
$ws.Range("N10").Value = 211.77486
$ws.Range("O10").Value = 0.3582195226742806
$ws.Range("P10").Value = 0.3859834752997077
$ws.Range("Q10").Value = 225.77402534544
$ws.Range("R10").Value = 2031.96622810896
$ws.Range("S10").Value = 0.05980484478133079
$ws.Range("T10").Value = 0.07258798185986376

$ws.Range("I11").Value = 0.1669502665149541
$ws.Range("J11").Value = 0.1880598173367416
$ws.Range("M11").Value = 7.337415
$ws.Range("N11").Value = 14.67483
$ws.Range("O11").Value = 0.03723395636710288
$ws.Range("P11").Value = 0.0267465263952124
$ws.Range("Q11").Value = 23.46734244348
$ws.Range("R11").Value = 140.80405466088
$ws.Range("S11").Value = 0.006216218938893997
$ws.Range("T11").Value = 0.005029946868275983

$ws.Range("I12").Value = 0.1669502665149541
$ws.Range("J12").Value = 0.1880598173367416
$ws.Range("M12").Value = 15.36873766666667
$ws.Range("N12").Value = 46.106213
$ws.Range("O12").Value = 0.07798917025929666
$ws.Range("P12").Value = 0.0840337532351506
$ws.Range("Q12").Value = 49.154018104152
$ws.Range("R12").Value = 442.386162937368
$ws.Range("S12").Value = 0.01302031276006971
$ws.Range("T12").Value = 0.01580337228352324

$ws.Range("I13").Value = 0.1669502665149541
$ws.Range("J13").Value = 0.1880598173367416
$ws.Range("M13").Value = 63.21821133333334
$ws.Range("N13").Value = 189.654634
$ws.Range("O13").Value = 0.3208029152489838
$ws.Range("P13").Value = 0.3456668782027013
$ws.Range("Q13").Value = 202.191563925936
$ws.Range("R13").Value = 1819.724075333424
$ws.Range("S13").Value = 0.05355813219959207
$ws.Range("T13").Value = 0.06500604997416172

$ws.Range("G14").Value = 2.617047
$ws.Range("H14").Value = 7.851141
$ws.Range("I14").Value = 0.1366085279147753
$ws.Range("J14").Value = 0.1538816040404024
$ws.Range("M14").Value = 35.18694
$ws.Range("N14").Value = 70.37388
$ws.Range("O14").Value = 0.1785572969025014
$ws.Range("P14").Value = 0.1282643028201015
$ws.Range("Q14").Value = 92.08587576618
$ws.Range("R14").Value = 552.51525459708
$ws.Range("S14").Value = 0.02439244947829219
$ws.Range("T14").Value = 0.01973751665908113

$ws.Range("G15").Value = 2.617047
$ws.Range("H15").Value = 7.851141
$ws.Range("I15").Value = 0.1366085279147753
$ws.Range("J15").Value = 0.1538816040404024
$ws.Range("O15").Value = 0.02719713854783473
$ws.Range("P15").Value = 0.02930506404712654
$ws.Range("Q15").Value = 14.026154993145
$ws.Range("R15").Value = 126.235394938305
$ws.Range("S15").Value = 0.003715361060513893
$ws.Range("T15").Value = 0.004509510262078561

$ws.Range("G16").Value = 2.617047
$ws.Range("H16").Value = 7.851141
$ws.Range("I16").Value = 0.1366085279147753
$ws.Range("J16").Value = 0.1538816040404024
$ws.Range("M16").Value = 70.59161999999999
$ws.Range("N16").Value = 211.77486
$ws.Range("O16").Value = 0.3582195226742806
$ws.Range("P16").Value = 0.3859834752997077
$ws.Range("Q16").Value = 184.74158734614
$ws.Range("R16").Value = 1662.67428611526
$ws.Range("S16").Value = 0.04893584166286697
$ws.Range("T16").Value = 0.05939575631220808

$ws.Range("G17").Value = 2.617047
$ws.Range("H17").Value = 7.851141
$ws.Range("I17").Value = 0.1366085279147753
$ws.Range("J17").Value = 0.1538816040404024
$ws.Range("M17").Value = 7.337415
$ws.Range("N17").Value = 14.67483
$ws.Range("O17").Value = 0.03723395636710288
$ws.Range("P17").Value = 0.0267465263952124
$ws.Range("Q17").Value = 19.202359913505
$ws.Range("R17").Value = 115.21415948103
$ws.Range("S17").Value = 0.005086475967752901
$ws.Range("T17").Value = 0.004115798384204247

$ws.Range("G18").Value = 2.617047
$ws.Range("H18").Value = 7.851141
$ws.Range("I18").Value = 0.1366085279147753
$ws.Range("J18").Value = 0.1538816040404024
$ws.Range("M18").Value = 15.36873766666667
$ws.Range("N18").Value = 46.106213
$ws.Range("O18").Value = 0.07798917025929666
$ws.Range("P18").Value = 0.0840337532351506
$ws.Range("Q18").Value = 40.220708804337
$ws.Range("R18").Value = 361.986379239033
$ws.Range("S18").Value = 0.0106539857424173
$ws.Range("T18").Value = 0.01293124874136033

$ws.Range("G19").Value = 2.617047
$ws.Range("H19").Value = 7.851141
$ws.Range("I19").Value = 0.1366085279147753
$ws.Range("J19").Value = 0.1538816040404024
$ws.Range("M19").Value = 63.21821133333334
$ws.Range("N19").Value = 189.654634
$ws.Range("O19").Value = 0.3208029152489838
$ws.Range("P19").Value = 0.3456668782027013
$ws.Range("Q19").Value = 165.445030315266
$ws.Range("R19").Value = 1489.005272837394
$ws.Range("S19").Value = 0.04382441400293211
$ws.Range("T19").Value = 0.0531917736814701

$ws.Range("G20").Value = 4.784714
$ws.Range("H20").Value = 9.569428
$ws.Range("I20").Value = 0.2497596474320929
$ws.Range("J20").Value = 0.1875598630045162
$ws.Range("M20").Value = 35.18694
$ws.Range("N20").Value = 70.37388
$ws.Range("O20").Value = 0.1785572969025014
$ws.Range("P20").Value = 0.1282643028201015
$ws.Range("Q20").Value = 168.35944443516
$ws.Range("R20").Value = 673.4377777406401
$ws.Range("S20").Value = 0.04459640752079627
$ws.Range("T20").Value = 0.02405723506530801

$ws.Range("G21").Value = 4.784714
$ws.Range("H21").Value = 9.569428
$ws.Range("I21").Value = 0.2497596474320929
$ws.Range("J21").Value = 0.1875598630045162
$ws.Range("O21").Value = 0.02719713854783473
$ws.Range("P21").Value = 0.02930506404712654
$ws.Range("Q21").Value = 25.64384214799
$ws.Range("R21").Value = 153.86305288794
$ws.Range("S21").Value = 0.006792747734868984
$ws.Range("T21").Value = 0.005496453798017628

$ws.Range("G22").Value = 4.784714
$ws.Range("H22").Value = 9.569428
$ws.Range("I22").Value = 0.2497596474320929
$ws.Range("J22").Value = 0.1875598630045162
$ws.Range("M22").Value = 70.59161999999999
$ws.Range("N22").Value = 211.77486
$ws.Range("O22").Value = 0.3582195226742806
$ws.Range("P22").Value = 0.3859834752997077
$ws.Range("Q22").Value = 337.76071249668
$ws.Range("R22").Value = 2026.56427498008
$ws.Range("S22").Value = 0.08946878168642095
$ws.Range("T22").Value = 0.07239500774922024

$ws.Range("G23").Value = 4.784714
$ws.Range("H23").Value = 9.569428
$ws.Range("I23").Value = 0.2497596474320929
$ws.Range("J23").Value = 0.1875598630045162
$ws.Range("M23").Value = 7.337415
$ws.Range("N23").Value = 14.67483
$ws.Range("O23").Value = 0.03723395636710288
$ws.Range("P23").Value = 0.0267465263952124
$ws.Range("Q23").Value = 35.10743227431
$ws.Range("R23").Value = 140.42972909724
$ws.Range("S23").Value = 0.009299539814749546
$ws.Range("T23").Value = 0.005016574826532714

$ws.Range("G24").Value = 4.784714
$ws.Range("H24").Value = 9.569428
$ws.Range("I24").Value = 0.2497596474320929
$ws.Range("J24").Value = 0.1875598630045162
$ws.Range("M24").Value = 15.36873766666667
$ws.Range("N24").Value = 46.106213
$ws.Range("O24").Value = 0.07798917025929666
$ws.Range("P24").Value = 0.0840337532351506
$ws.Range("Q24").Value = 73.53501427602734
$ws.Range("R24").Value = 441.210085656164
$ws.Range("S24").Value = 0.0194785476674834
$ws.Range("T24").Value = 0.01576135924454016

$ws.Range("G25").Value = 4.784714
$ws.Range("H25").Value = 9.569428
$ws.Range("I25").Value = 0.2497596474320929
$ws.Range("J25").Value = 0.1875598630045162
$ws.Range("M25").Value = 63.21821133333334
$ws.Range("N25").Value = 189.654634
$ws.Range("O25").Value = 0.3208029152489838
$ws.Range("P25").Value = 0.3456668782027013
$ws.Range("Q25").Value = 302.4810608215587
$ws.Range("R25").Value = 1814.886364929352
$ws.Range("S25").Value = 0.08012362300777377
$ws.Range("T25").Value = 0.06483323232089744

$ws.Range("G26").Value = 4.899255333333334
$ws.Range("H26").Value = 14.697766
$ws.Range("I26").Value = 0.2557386470190557
$ws.Range("J26").Value = 0.2880747916628283
$ws.Range("M26").Value = 35.18694
$ws.Range("N26").Value = 70.37388
$ws.Range("O26").Value = 0.1785572969025014
$ws.Range("P26").Value = 0.1282643028201015
$ws.Range("Q26").Value = 172.38980345868
$ws.Range("R26").Value = 1034.33882075208
$ws.Range("S26").Value = 0.04566400152522552
$ws.Range("T26").Value = 0.03694971231267866

$ws.Range("G27").Value = 4.899255333333334
$ws.Range("H27").Value = 14.697766
$ws.Range("I27").Value = 0.2557386470190557
$ws.Range("J27").Value = 0.2880747916628283
$ws.Range("O27").Value = 0.02719713854783473
$ws.Range("P27").Value = 0.02930506404712654
$ws.Range("Q27").Value = 26.25773043293667
$ws.Range("R27").Value = 236.31957389643
$ws.Range("S27").Value = 0.006955359415013058
$ws.Range("T27").Value = 0.008442050220041819

$ws.Range("G28").Value = 4.899255333333334
$ws.Range("H28").Value = 14.697766
$ws.Range("I28").Value = 0.2557386470190557
$ws.Range("J28").Value = 0.2880747916628283
$ws.Range("M28").Value = 70.59161999999999
$ws.Range("N28").Value = 211.77486
$ws.Range("O28").Value = 0.3582195226742806
$ws.Range("P28").Value = 0.3859834752997077
$ws.Range("Q28").Value = 345.84637077364
$ws.Range("R28").Value = 3112.61733696276
$ws.Range("S28").Value = 0.09161057606453248
$ws.Range("T28").Value = 0.1111921092322577

$ws.Range("G29").Value = 4.899255333333334
$ws.Range("H29").Value = 14.697766
$ws.Range("I29").Value = 0.2557386470190557
$ws.Range("J29").Value = 0.2880747916628283
$ws.Range("M29").Value = 7.337415
$ws.Range("N29").Value = 14.67483
$ws.Range("O29").Value = 0.03723395636710288
$ws.Range("P29").Value = 0.0267465263952124
$ws.Range("Q29").Value = 35.94786957163
$ws.Range("R29").Value = 215.68721742978
$ws.Range("S29").Value = 0.009522161624489444
$ws.Range("T29").Value = 0.00770500001900515

$ws.Range("G30").Value = 4.899255333333334
$ws.Range("H30").Value = 14.697766
$ws.Range("I30").Value = 0.2557386470190557
$ws.Range("J30").Value = 0.2880747916628283
$ws.Range("M30").Value = 15.36873766666667
$ws.Range("N30").Value = 46.106213
$ws.Range("O30").Value = 0.07798917025929666
$ws.Range("P30").Value = 0.0840337532351506
$ws.Range("Q30").Value = 75.29536998001755
$ws.Range("R30").Value = 677.658329820158
$ws.Range("S30").Value = 0.0199448448842513
$ws.Range("T30").Value = 0.02420800595586153

$ws.Range("G31").Value = 4.899255333333334
$ws.Range("H31").Value = 14.697766
$ws.Range("I31").Value = 0.2557386470190557
$ws.Range("J31").Value = 0.2880747916628283
$ws.Range("M31").Value = 63.21821133333334
$ws.Range("N31").Value = 189.654634
$ws.Range("O31").Value = 0.3208029152489838
$ws.Range("P31").Value = 0.3456668782027013
$ws.Range("Q31").Value = 309.7221590386271
$ws.Range("R31").Value = 2787.499431347644
$ws.Range("S31").Value = 0.08204170350554391
$ws.Range("T31").Value = 0.09957791392298342

$ws.Range("G32").Value = 1.991497666666667
$ws.Range("H32").Value = 5.974493
$ws.Range("I32").Value = 0.1039551695437809
$ws.Range("J32").Value = 0.117099484796943
$ws.Range("M32").Value = 35.18694
$ws.Range("N32").Value = 70.37388
$ws.Range("O32").Value = 0.1785572969025014
$ws.Range("P32").Value = 0.1282643028201015
$ws.Range("Q32").Value = 70.07470890713999
$ws.Range("R32").Value = 420.44825344284
$ws.Range("S32").Value = 0.01856195407277876
$ws.Range("T32").Value = 0.01501968377807297

$ws.Range("G33").Value = 1.991497666666667
$ws.Range("H33").Value = 5.974493
$ws.Range("I33").Value = 0.1039551695437809
$ws.Range("J33").Value = 0.117099484796943
$ws.Range("O33").Value = 0.02719713854783473
$ws.Range("P33").Value = 0.02930506404712654
$ws.Range("Q33").Value = 10.67350144691833
$ws.Range("R33").Value = 96.06151302226499
$ws.Range("S33").Value = 0.002827283148845859
$ws.Range("T33").Value = 0.003431607901859937

$ws.Range("G34").Value = 1.991497666666667
$ws.Range("H34").Value = 5.974493
$ws.Range("I34").Value = 0.1039551695437809
$ws.Range("J34").Value = 0.117099484796943
$ws.Range("M34").Value = 70.59161999999999
$ws.Range("N34").Value = 211.77486
$ws.Range("O34").Value = 0.3582195226742806
$ws.Range("P34").Value = 0.3859834752997077
$ws.Range("Q34").Value = 140.58304651622
$ws.Range("R34").Value = 1265.24741864598
$ws.Range("S34").Value = 0.03723877121349713
$ws.Range("T34").Value = 0.04519846609772936

$ws.Range("G35").Value = 1.991497666666667
$ws.Range("H35").Value = 5.974493
$ws.Range("I35").Value = 0.1039551695437809
$ws.Range("J35").Value = 0.117099484796943
$ws.Range("M35").Value = 7.337415
$ws.Range("N35").Value = 14.67483
$ws.Range("O35").Value = 0.03723395636710288
$ws.Range("P35").Value = 0.0267465263952124
$ws.Range("Q35").Value = 14.612444851865
$ws.Range("R35").Value = 87.67466911119
$ws.Range("S35").Value = 0.003870662246927922
$ws.Range("T35").Value = 0.00313200446098721

$ws.Range("G36").Value = 1.991497666666667
$ws.Range("H36").Value = 5.974493
$ws.Range("I36").Value = 0.1039551695437809
$ws.Range("J36").Value = 0.117099484796943
$ws.Range("M36").Value = 15.36873766666667
$ws.Range("N36").Value = 46.106213
$ws.Range("O36").Value = 0.07798917025929666
$ws.Range("P36").Value = 0.0840337532351506
$ws.Range("Q36").Value = 30.60680520277878
$ws.Range("R36").Value = 275.461246825009
$ws.Range("S36").Value = 0.008107377416883983
$ws.Range("T36").Value = 0.009840309209389579

$ws.Range("G37").Value = 1.991497666666667
$ws.Range("H37").Value = 5.974493
$ws.Range("I37").Value = 0.1039551695437809
$ws.Range("J37").Value = 0.117099484796943
$ws.Range("M37").Value = 63.21821133333334
$ws.Range("N37").Value = 189.654634
$ws.Range("O37").Value = 0.3208029152489838
$ws.Range("P37").Value = 0.3456668782027013
$ws.Range("Q37").Value = 125.8989203611736
$ws.Range("R37").Value = 1133.090283250562
$ws.Range("S37").Value = 0.0333491214448473
$ws.Range("T37").Value = 0.04047741334890398
